$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6855646666666667
$ws.Range("H2").Value = 2.056694
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2514313333333333
$ws.Range("N2").Value = 0.754294
$ws.Range("O2").Value = 0.01999844741031966
$ws.Range("P2").Value = 0.01999844741031965
$ws.Range("Q2").Value = 0.1723724382262222
$ws.Range("R2").Value = 1.551351944036
$ws.Range("S2").Value = 0.01999844741031966
$ws.Range("T2").Value = 0.01999844741031965

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6855646666666667
$ws.Range("H3").Value = 2.056694
$ws.Range("O3").Value = 0.790117395505359
$ws.Range("P3").Value = 0.7901173955053589
$ws.Range("Q3").Value = 6.810251773741777
$ws.Range("R3").Value = 61.292265963676
$ws.Range("S3").Value = 0.790117395505359
$ws.Range("T3").Value = 0.7901173955053589

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6855646666666667
$ws.Range("H4").Value = 2.056694
$ws.Range("M4").Value = 2.387326666666667
$ws.Range("N4").Value = 7.16198
$ws.Range("O4").Value = 0.1898841570843214
$ws.Range("P4").Value = 0.1898841570843214
$ws.Range("Q4").Value = 1.636666810457778
$ws.Range("R4").Value = 14.73000129412
$ws.Range("S4").Value = 0.1898841570843214
$ws.Range("T4").Value = 0.1898841570843214
